$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two FedEx tracking numbers on rows 2 and 3 (column P).
# The values are numeric-looking strings that must remain stored as
# shared-string text (matching the original cell type/format), so we
# write them via a text formula and then convert the formula result to
# a plain value with Paste Special (values only). This avoids Excel's
# "number stored as text" quote-prefix/number-format side effects that
# a plain Value assignment of a digit-only string would trigger.

$ws.Range("P2").Formula = '="320018253468"'
$ws.Range("P2").Copy()
$ws.Range("P2").PasteSpecial(-4163)

$ws.Range("P3").Formula = '="320018253479"'
$ws.Range("P3").Copy()
$ws.Range("P3").PasteSpecial(-4163)

$excel.CutCopyMode = 0
